# Apply the latest coinranking.com snapshot values to the cryptos sheet.
# Column D ("Price") values are digit/period strings that must stay text
# (e.g. "28.434.15", "0.9996"); Column B/C/E are plain text already, so a
# plain .Value assignment is enough for those.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force text storage so numeric-looking strings (prices such as
    # "0.9996" or "1.811.71") are not reinterpreted as numbers/dates.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue 'D2' '28.434.15'
$ws.Range('E2').Value = '  -0.18%  '
Set-TextValue 'D3' '1.811.71'
$ws.Range('E3').Value = '  -0.80%  '
Set-TextValue 'D4' '0.9996'
$ws.Range('E4').Value = '  -0.44%  '
Set-TextValue 'D5' '313.02'
$ws.Range('E5').Value = '  -1.08%  '
Set-TextValue 'D6' '0.9998'
$ws.Range('E6').Value = '  -0.35%  '
Set-TextValue 'D7' '0.5148'
$ws.Range('E7').Value = '  -0.48%  '
Set-TextValue 'D8' '0.4001'
$ws.Range('E8').Value = '  +3.57%  '
Set-TextValue 'D9' '0.07874'
$ws.Range('E9').Value = '  -5.00%  '
Set-TextValue 'D10' '1.116'
$ws.Range('E10').Value = '  -0.70%  '
Set-TextValue 'D11' '41.02'
$ws.Range('E11').Value = '  -2.18%  '
Set-TextValue 'D12' '6.384'
$ws.Range('E12').Value = '  -0.04%  '
Set-TextValue 'D13' '0.9997'
$ws.Range('E13').Value = '  -0.37%  '
Set-TextValue 'D14' '20.43'
$ws.Range('E14').Value = '  -3.72%  '
$ws.Range('E15').Value = '  -2.22%  '
Set-TextValue 'D16' '1.809.58'
$ws.Range('E16').Value = '  -1.02%  '
Set-TextValue 'D17' '92.90'
$ws.Range('E17').Value = '  -1.13%  '
Set-TextValue 'D18' '0.00001084'
$ws.Range('E18').Value = '  -3.47%  '
$ws.Range('E19').Value = '  -0.90%  '
Set-TextValue 'D20' '0.9995'
$ws.Range('E20').Value = '  -0.36%  '
Set-TextValue 'D21' '17.34'
$ws.Range('E21').Value = '  -2.64%  '
Set-TextValue 'D22' '6.020'
$ws.Range('E22').Value = '  -0.69%  '
Set-TextValue 'D23' '28.475.54'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('E24').Value = '  -2.80%  '
Set-TextValue 'D25' '2.234'
$ws.Range('E25').Value = '  -0.48%  '
Set-TextValue 'D26' '161.77'
$ws.Range('E26').Value = '  +1.16%  '
Set-TextValue 'D27' '20.55'
$ws.Range('E27').Value = '  -2.55%  '
Set-TextValue 'D28' '2.020.25'
$ws.Range('E28').Value = '  -0.86%  '
Set-TextValue 'D29' '2.410'
$ws.Range('E29').Value = '  -0.31%  '
Set-TextValue 'D30' '128.34'
$ws.Range('E30').Value = '  +1.99%  '
Set-TextValue 'D31' '0.1091'
$ws.Range('E31').Value = '  -0.45%  '
Set-TextValue 'D32' '1.071'
$ws.Range('E32').Value = '  -2.48%  '
$ws.Range('E33').Value = '  -0.47%  '
Set-TextValue 'D34' '5.589'
$ws.Range('E34').Value = '  -2.50%  '
Set-TextValue 'D35' '0.07259'
$ws.Range('E35').Value = '  -5.04%  '
Set-TextValue 'D36' '9.287'
$ws.Range('E36').Value = '  +5.70%  '
$ws.Range('E37').Value = '  -1.25%  '
Set-TextValue 'D38' '0.2176'
$ws.Range('E38').Value = '  -2.64%  '
Set-TextValue 'D39' '11.70'
$ws.Range('E39').Value = '  -2.89%  '
Set-TextValue 'D40' '5.063'
$ws.Range('E40').Value = '  -3.78%  '
Set-TextValue 'D41' '0.6210'
$ws.Range('E41').Value = '  -3.37%  '
Set-TextValue 'D42' '0.9993'
$ws.Range('E42').Value = '  -0.40%  '
Set-TextValue 'D43' '1.161'
$ws.Range('E43').Value = '  -2.64%  '
Set-TextValue 'D44' '13.22'
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D45' '0.6005'
$ws.Range('E45').Value = '  -3.34%  '
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D46' '1.315'
$ws.Range('E46').Value = '  -6.05%  '
Set-TextValue 'D47' '3.733'
$ws.Range('E47').Value = '  -1.69%  '
Set-TextValue 'D48' '125.99'
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('E49').Value = '  +1.57%  '
Set-TextValue 'D50' '1.933'
$ws.Range('E50').Value = '  -3.57%  '
Set-TextValue 'D51' '0.06849'
